# "add tabel format baru" — renumber the Kecamatan-level health tables
# (Tabel 4.2.3/4.2.4/4.2.5 -> 4.2.5/4.2.6/4.2.7) and bump the reporting
# year in the table captions/titles from 2020 to 2021.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bab 4")

# --- Table 2 header (H1): "Tabel 4.2.3" -> "Tabel 4.2.5" ------------------
$ws.Range("H1").Value = "Tabel 4.2.5"

# --- Table titles: bump the year 2020 -> 2021 ------------------------------
$ws.Range("B1").Value = "Jumlah Tenaga Kesehatan Menurut Kelurahan/Desa in Kecamatan Toari. 2021"
$ws.Range("I1").Value = "Banyaknya Bayi yang Diimunisasi Menurut Jenis dan Desa/Kelurahan di Kecamatan Toari, 2021"
$ws.Range("Q1").Value = "Banyaknya Ibu Melahirkan dan Kelahiran Ditolong Tenaga Kesehatan Menurut Desa/Kelurahan di Kecamatan Toari, 2021"
$ws.Range("X1").Value = "Banyaknya Pasangan Usia Subur dan Peserta KB Menurut Desa/Kelurahan di Kecamatan Toari, 2021"

# --- Table 3 header (P1): rich text "Tabel" + " 4.2.4." -> " 4.2.6." ------
$p1 = $ws.Range("P1")
$p1Len = $p1.Characters().Text.Length
$p1Suffix = $p1.Characters(6, $p1Len - 5)
$p1Suffix.Text = " 4.2.6."
$p1Suffix.Font.Underline = $false
$p1Suffix.Font.Size = 9

# --- Table 4 header (W1): rich text "Tabel" + " 4.2.5." -> " 4.2.7." ------
$w1 = $ws.Range("W1")
$w1Len = $w1.Characters().Text.Length
$w1Suffix = $w1.Characters(6, $w1Len - 5)
$w1Suffix.Text = " 4.2.7."
$w1Suffix.Font.Underline = $false
$w1Suffix.Font.Size = 9

# --- Reset the sheet view: clear the scrolled/selected state left over
#     from editing (topLeftCell=A10 / selection on J21) back to A1. --------
$ws.Range("A1").Select()
